# Add a new row (row 8) to the active sheet, mirroring row 7's data style
# for column A, and filling B:G with text values, matching the author's
# commit that appended a new "ss88" record to the students table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last data row's first cell (A7) onto the new
# row's first cell (A8) so the new id cell keeps the same style (s="1").
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row's values.
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "ss88"
$ws.Range("C8").Value = "ss88"
$ws.Range("D8").Value = "ss88"
$ws.Range("E8").Value = "ss88ss88"
$ws.Range("F8").Value = "ss88"
$ws.Range("G8").Value = "ss88"
